$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$suffix = ".age_premiere_conso"

# Columns B through P on row 1 all currently end with ".deja" and need
# ".age_premiere_conso" appended right after it.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")

foreach ($col in $cols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = [string]$cell.Value2 + $suffix
}
